$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 307, pushing the existing 307:378 block down to 310:381.
$ws.Rows("307:309").Insert()

# Row 307 (new) — Acelga, Extra
$ws.Range("A307").Value = 9
$ws.Range("B307").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C307").Value = "Metropolitana"
$ws.Range("D307").Value = 44511
$ws.Range("E307").Value = 13
$ws.Range("F307").Value = 100112009
$ws.Range("G307").Value = "Acelga"
$ws.Range("H307").Value = "Sin especificar"
$ws.Range("I307").Value = "Extra"
$ws.Range("J307").Value = 43
$ws.Range("K307").Value = 12000
$ws.Range("L307").Value = 12000
$ws.Range("M307").Value = 12000
$ws.Range("N307").Value = "`$/docena de atados"
$ws.Range("O307").Value = "Región Metropolitana"
$ws.Range("P307").Value = 4000
$ws.Range("Q307").Value = 3
$ws.Range("R307").Value = "Hortaliza"

# Row 308 (new) — Acelga, Primera
$ws.Range("A308").Value = 9
$ws.Range("B308").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C308").Value = "Metropolitana"
$ws.Range("D308").Value = 44511
$ws.Range("E308").Value = 13
$ws.Range("F308").Value = 100112009
$ws.Range("G308").Value = "Acelga"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 61
$ws.Range("K308").Value = 10000
$ws.Range("L308").Value = 11000
$ws.Range("M308").Value = 10492
$ws.Range("N308").Value = "`$/docena de atados"
$ws.Range("O308").Value = "Región Metropolitana"
$ws.Range("P308").Value = 3497
$ws.Range("Q308").Value = 3
$ws.Range("R308").Value = "Hortaliza"

# Row 309 (new) — Acelga, Segunda
$ws.Range("A309").Value = 9
$ws.Range("B309").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C309").Value = "Metropolitana"
$ws.Range("D309").Value = 44511
$ws.Range("E309").Value = 13
$ws.Range("F309").Value = 100112009
$ws.Range("G309").Value = "Acelga"
$ws.Range("H309").Value = "Sin especificar"
$ws.Range("I309").Value = "Segunda"
$ws.Range("J309").Value = 34
$ws.Range("K309").Value = 8000
$ws.Range("L309").Value = 9000
$ws.Range("M309").Value = 8500
$ws.Range("N309").Value = "`$/docena de atados"
$ws.Range("O309").Value = "Región Metropolitana"
$ws.Range("P309").Value = 2833
$ws.Range("Q309").Value = 3
$ws.Range("R309").Value = "Hortaliza"
